{"js": "const results = context.document.body.search(\"01 November 2022\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n\n  // Replace the whole text with the first segment \"01 \" (keeps original run formatting).\n  const firstRange = range.insertText(\"01 \", \"Replace\");\n  await context.sync();\n\n  // Insert \"OCTOBER\" as its own run right after \"01 \".\n  const octRange = firstRange.insertText(\"OCTOBER\", \"After\");\n  await context.sync();\n\n  // Insert \" 2022\" as its own run right after \"OCTOBER\".\n  octRange.insertText(\" 2022\", \"After\");\n  await context.sync();\n}\n", "ps1": "# Ideation phase / LITERATURE SURVEY.docx\n# \"Date\" table cell: \"01 November 2022\" -> \"01 \" + \"OCTOBER\" + \" 2022\"\n# (three separate runs, matching the author's in-place retyping of the month).\n\n$d = $word.ActiveDocument\n\n$target = $d.Content\n$target.Find.ClearFormatting()\n$target.Find.Text = \"01 November 2022\"\n$target.Find.Execute() | Out-Null\n\n# Rebuild the paragraph with the original run (\"01 \") kept as-is (same rsid\n# attribute, just shorter text) followed by two brand-new runs (\"OCTOBER\" and\n# \" 2022\") that carry the identical run formatting. Using InsertXML (instead\n# of Range.Text / InsertAfter) keeps the three runs distinct instead of\n# Word's usual \"merge adjacent runs with identical formatting\" behaviour.\n$xml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w14:paraId=\"22D9CE15\" w14:textId=\"77777777\" w:rsidR=\"00970245\" w:rsidRPr=\"00970245\" w:rsidRDefault=\"00970245\" w:rsidP=\"00970245\"><w:pPr><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r w:rsidRPr=\"00970245\"><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">01 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>OCTOBER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> 2022</w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($xml)\n"}
